# Applies the "Kaos/REKAPITULASI" revision:
#  - header label BANYAK -> KIRIM (D1)
#  - removes the old "row 13" transaction (4 / WEARPACK 70 / F13 formula)
#  - adds a new "TOPI" line (row 12) that finishes transaction #3
#  - rebuilds transaction #4 (rows 14-16: WEARPACK / KAOS / TOPI)
#  - adds a brand new transaction #5 (rows 18-20: WEARPACK / KAOS / TOPI)
#  - adds a standalone total formula in F23
#  - keeps the view pointed at the new bottom of the sheet

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# grab the existing number-format strings so re-used styles map onto the
# same style indexes the workbook already has (comma style, two date styles)
$commaFmt = $ws1.Range("F15").NumberFormat
$date1Fmt = $ws1.Range("B10").NumberFormat
$date2Fmt = $ws1.Range("B13").NumberFormat

# ---- header row -----------------------------------------------------
$ws1.Range("D1").Value = "KIRIM"

# ---- drop the old row 13 (4 / 03-Apr / WEARPACK / 70 / =F11+D13) ----
$ws1.Rows(13).Clear()

# ---- new row 12: closes out transaction #3 with a TOPI line --------
$ws1.Range("C12").Value = "TOPI"
$ws1.Range("F12").Value = 10
$ws1.Range("F12").NumberFormat = $commaFmt

# ---- row 14: transaction #4 header + WEARPACK line ------------------
$ws1.Range("A14").Value = 4
$ws1.Range("B14").Value = 44289
$ws1.Range("B14").NumberFormat = $date2Fmt
$ws1.Range("C14").Value = "WEARPACK"
$ws1.Range("D14").Value = 70
$ws1.Range("F14").Formula = "=F11+D14"
$ws1.Range("F14").NumberFormat = $commaFmt

# ---- row 15: KAOS line (already existed, keep as-is) ----------------
$ws1.Range("C15").Value = "KAOS"
$ws1.Range("D15").Value = 70
$ws1.Range("F15").Value = 110
$ws1.Range("F15").NumberFormat = $commaFmt

# ---- row 16: new TOPI line closing transaction #4 --------------------
$ws1.Range("C16").Value = "TOPI"
$ws1.Range("D16").Value = 70
$ws1.Range("F16").Value = 80
$ws1.Range("F16").NumberFormat = $commaFmt

# ---- row 18-20: brand new transaction #5 -----------------------------
$ws1.Range("A18").Value = 5
$ws1.Range("B18").Value = 44298
$ws1.Range("B18").NumberFormat = $date2Fmt
$ws1.Range("C18").Value = "WEARPACK"
$ws1.Range("D18").Value = 120
$ws1.Range("F18").Value = 230
$ws1.Range("F18").NumberFormat = $commaFmt

$ws1.Range("C19").Value = "KAOS"
$ws1.Range("D19").Value = 120
$ws1.Range("F19").Value = 230
$ws1.Range("F19").NumberFormat = $commaFmt

$ws1.Range("C20").Value = "TOPI"
$ws1.Range("D20").Value = 120
$ws1.Range("F20").Formula = "=F16+D20"
$ws1.Range("F20").NumberFormat = $commaFmt

# ---- row 23: standalone grand-total formula --------------------------
$ws1.Range("F23").Formula = "=350*350000"
$ws1.Range("F23").NumberFormat = $commaFmt

# ---- refresh the view so it matches the newly extended sheet ---------
$ws1.Activate() | Out-Null
$ws1.Range("F23").Select() | Out-Null
$excel.ActiveWindow.Zoom = 216 | Out-Null
